$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 328 - this shifts the existing rows 328-400
# down to 329-401 (matching dimension change A1:R400 -> A1:R401).
$ws.Rows.Item(328).Insert()

# Populate the newly inserted row 328 with a fresh weekly data point for
# "Ajo" / "Chino" / "Primera" (same descriptive columns as the record that
# used to sit at row 328), with new date + price figures.
$ws.Range("A328").Value = 6
$ws.Range("B328").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C328").Value = "Metropolitana"
$ws.Range("D328").Value = 44508
$ws.Range("E328").Value = 13
$ws.Range("F328").Value = 100112003
$ws.Range("G328").Value = "Ajo"
$ws.Range("H328").Value = "Chino"
$ws.Range("I328").Value = "Primera"
$ws.Range("J328").Value = 246000
$ws.Range("K328").Value = 16500
$ws.Range("L328").Value = 17000
$ws.Range("M328").Value = 16553
$ws.Range("N328").Value = "$/caja 10 kilos"
$ws.Range("O328").Value = "China"
$ws.Range("P328").Value = 1655
$ws.Range("Q328").Value = 10
$ws.Range("R328").Value = "Hortaliza"
